$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the 2023 column (T): header (row 4) + value (row 5).
#    Clone the formatting of the neighbouring 2022 column (S) so the new
#    cells pick up the same number formats / fonts / borders, then set the
#    real values on top of the pasted format.
# ---------------------------------------------------------------------------
$ws.Range("S4").Copy() | Out-Null
$ws.Range("T4").PasteSpecial(-4122) | Out-Null
$ws.Range("S5").Copy() | Out-Null
$ws.Range("T5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 105.59374642341281

# ---------------------------------------------------------------------------
# 2) Footnote row (row 6): Kyrgyz / Russian / English text with a
#    superscript "1" marker at the start of each cell, in columns A/B/C.
# ---------------------------------------------------------------------------
$footnoteKg = [string][char]0x0031 + " 2020-жылдан баштап маалыматтар, 2008 жылдагы Улуттук Эсептер Тутумунун эл аралык стандарттарына ылайык эсептелген "
$footnoteRu = [string][char]0x0031 + " Данные с 2020 года рассчитаны по международному стандарту Системы Национальных Счетов 2008 года"
$footnoteEn = [string][char]0x0031 + " Data from 2020 are calculated according to the international standard of the System of National Accounts 2008"

function Set-Footnote($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.Value = $text
    $len = $text.Length

    # Leading "1" -> superscript, same italic 8pt Times New Roman run font
    $sup = $rng.Characters(1, 1)
    $sup.Font.Superscript = $true
    $sup.Font.Italic = $true
    $sup.Font.Size = 8
    $sup.Font.Name = "Times New Roman"
    $sup.Font.ColorIndex = -4105

    # Remainder of the text -> plain italic 8pt Times New Roman run font
    $body = $rng.Characters(2, $len - 1)
    $body.Font.Italic = $true
    $body.Font.Size = 8
    $body.Font.Name = "Times New Roman"
    $body.Font.ColorIndex = -4105

    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
    $rng.NumberFormat = "General"
}

Set-Footnote "A6" $footnoteKg
Set-Footnote "B6" $footnoteRu
Set-Footnote "C6" $footnoteEn

$ws.Range("A6:M6").RowHeight = 36.75

Write-Output "done"
